$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''69.127.67'
$ws.Range('E2').Value = '  -3.52%  '

$ws.Range('D3').Value = '''3.519.98'
$ws.Range('E3').Value = '  -4.60%  '

$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').Value = '''581.51'
$ws.Range('E5').Value = '  -1.26%  '

$ws.Range('D6').Value = '''175.66'
$ws.Range('E6').Value = '  -2.22%  '

$ws.Range('D7').Value = '''0.624'
$ws.Range('E7').Value = '  +0.15%  '

$ws.Range('D8').Value = '''3.513.34'
$ws.Range('E8').Value = '  -4.57%  '

$ws.Range('D10').Value = '''0.190'
$ws.Range('E10').Value = '  -5.62%  '

$ws.Range('D11').Value = '''6.74'
$ws.Range('E11').Value = '  +7.74%  '

$ws.Range('E12').Value = '  -1.92%  '

$ws.Range('D13').Value = '''47.43'
$ws.Range('E13').Value = '  -4.92%  '

$ws.Range('E14').Value = '  -2.85%  '

$ws.Range('D15').Value = '''673.33'
$ws.Range('E15').Value = '  -1.34%  '

$ws.Range('D16').Value = '''4.085.04'
$ws.Range('E16').Value = '  -4.70%  '

$ws.Range('E17').Value = '  -1.82%  '

$ws.Range('D18').Value = '''3.518.74'
$ws.Range('E18').Value = '  -4.65%  '

$ws.Range('D19').Value = '''69.094.21'
$ws.Range('E19').Value = '  -3.75%  '

$ws.Range('E20').Value = '  -1.49%  '

$ws.Range('D21').Value = '''17.63'
$ws.Range('E21').Value = '  -2.42%  '

$ws.Range('D22').Value = '''11.30'
$ws.Range('E22').Value = '  -3.21%  '

$ws.Range('D23').Value = '''0.910'
$ws.Range('E23').Value = '  -3.15%  '

$ws.Range('D24').Value = '''16.33'
$ws.Range('E24').Value = '  -8.44%  '

$ws.Range('D25').Value = '''98.38'
$ws.Range('E25').Value = '  -5.23%  '

$ws.Range('E26').Value = '  -4.14%  '

$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').Value = '''5.84'
$ws.Range('E27').Value = '  +0.20%  '

$ws.Range('D28').Value = '''2.68'
$ws.Range('E28').Value = '  -5.89%  '

$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  +0.08%  '

$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').Value = '''9.53'
$ws.Range('E30').Value = '  -6.66%  '

$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '''33.04'
$ws.Range('E31').Value = '  -6.73%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '''8.79'
$ws.Range('E32').Value = '  -4.74%  '

$ws.Range('B33').Value = 'Stacks'
$ws.Range('C33').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D33').Value = '''3.23'
$ws.Range('E33').Value = '  -6.96%  '

$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '''7.36'
$ws.Range('E34').Value = '  -0.11%  '

$ws.Range('B35').Value = 'Mantle'
$ws.Range('C35').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D35').Value = '''1.36'
$ws.Range('E35').Value = '  -4.45%  '

$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').Value = '''579.94'
$ws.Range('E36').Value = '  +0.98%  '

$ws.Range('D37').Value = '''3.61'
$ws.Range('E37').Value = '  -14.11%  '

$ws.Range('B38').Value = 'Cosmos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D38').Value = '''10.97'
$ws.Range('E38').Value = '  -2.95%  '

$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '''0.106'
$ws.Range('E39').Value = '  -3.47%  '

$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = '''57.39'
$ws.Range('E40').Value = '  -3.47%  '

$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '''1.00'
$ws.Range('E41').Value = '  +0.06%  '

$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').Value = '''0.340'
$ws.Range('E42').Value = '  -3.40%  '

$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '''0.0441'
$ws.Range('E43').Value = '  -4.58%  '

$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').Value = '''0.137'
$ws.Range('E44').Value = '  -5.78%  '

$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '''3.428.07'
$ws.Range('E45').Value = '  -9.05%  '

$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = '''33.59'
$ws.Range('E46').Value = '  -5.10%  '

$ws.Range('B47').Value = 'PEPE'
$ws.Range('C47').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D47').Value = '''0.0₃0709'
$ws.Range('E47').Value = '  -8.57%  '

$ws.Range('B48').Value = 'ThetaToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D48').Value = '''2.91'
$ws.Range('E48').Value = '  +1.19%  '

$ws.Range('B49').Value = 'Fetch.AI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D49').Value = '''2.61'
$ws.Range('E49').Value = '  -6.52%  '

$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').Value = '''0.134'
$ws.Range('E50').Value = '  -0.15%  '

$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').Value = '''131.64'
$ws.Range('E51').Value = '  -1.96%  '
